# Auto-generated Excel COM-interop script
# Applies numeric updates to the Leve profit-tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 237.3
$ws.Range("I33").Value = 247.11111
$ws.Range("K33").Value = 247.11111
$ws.Range("M33").Value = -18.11111
$ws.Range("H46").Value = 100
$ws.Range("I46").Value = 100
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 300
$ws.Range("M46").Value = -181
$ws.Range("N46").Value = -538
$ws.Range("H55").Value = 1586
$ws.Range("I55").Value = 369.91666
$ws.Range("J55").Value = 3670.7144
$ws.Range("K55").Value = 369.91666
$ws.Range("L55").Value = 3670.7144
$ws.Range("M55").Value = -155.91666
$ws.Range("N55").Value = -4098.7144
$ws.Range("H59").Value = 2974.75
$ws.Range("J59").Value = 2966.3333
$ws.Range("L59").Value = 8898.999899999999
$ws.Range("N59").Value = -10012.9999
$ws.Range("H60").Value = 100
$ws.Range("I60").Value = 100
$ws.Range("J60").Value = 100
$ws.Range("K60").Value = 300
$ws.Range("L60").Value = 300
$ws.Range("M60").Value = 184
$ws.Range("N60").Value = -1268
$ws.Range("H61").Value = 2158.6667
$ws.Range("I61").Value = 2158.6667
$ws.Range("K61").Value = 6476.000100000001
$ws.Range("M61").Value = -6304.000100000001
$ws.Range("H132").Value = 2113.2727
$ws.Range("I132").Value = 2113.2727
$ws.Range("K132").Value = 6339.8181
$ws.Range("M132").Value = -3809.8181
$ws.Range("H137").Value = 5236.8887
$ws.Range("I137").Value = 4766.6665
$ws.Range("J137").Value = 5472
$ws.Range("K137").Value = 14299.9995
$ws.Range("L137").Value = 16416
$ws.Range("M137").Value = -11749.9995
$ws.Range("N137").Value = -21516

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -15976
$ws.Range("H61").Value = 3679.625
$ws.Range("I61").Value = 2978.739
$ws.Range("J61").Value = 19800
$ws.Range("K61").Value = 2978.739
$ws.Range("L61").Value = 19800
$ws.Range("M61").Value = -2766.739
$ws.Range("N61").Value = -20224
$ws.Range("H74").Value = 15153223
$ws.Range("I74").Value = 15874709
$ws.Range("J74").Value = 2014
$ws.Range("K74").Value = 15874709
$ws.Range("L74").Value = 2014
$ws.Range("M74").Value = -15873835
$ws.Range("N74").Value = -3762
$ws.Range("H77").Value = 15153223
$ws.Range("I77").Value = 15874709
$ws.Range("J77").Value = 2014
$ws.Range("K77").Value = 79373545
$ws.Range("L77").Value = 10070
$ws.Range("M77").Value = -79369177
$ws.Range("N77").Value = -18806
$ws.Range("H97").Value = 2006
$ws.Range("I97").Value = 2504.5
$ws.Range("J97").Value = 1009
$ws.Range("K97").Value = 2504.5
$ws.Range("L97").Value = 1009
$ws.Range("M97").Value = -2008.5
$ws.Range("N97").Value = -2001
$ws.Range("H122").Value = 1764.04
$ws.Range("I122").Value = 1420.875
$ws.Range("K122").Value = 4262.625
$ws.Range("M122").Value = -1812.625
$ws.Range("H127").Value = 53999.5
$ws.Range("J127").Value = 53999.5
$ws.Range("L127").Value = 53999.5
$ws.Range("N127").Value = -63919.5
$ws.Range("H128").Value = 99997.5
$ws.Range("J128").Value = 99997.5
$ws.Range("L128").Value = 99997.5
$ws.Range("N128").Value = -109957.5
$ws.Range("H132").Value = 4315.5713
$ws.Range("I132").Value = 2617
$ws.Range("J132").Value = 14507
$ws.Range("K132").Value = 7851
$ws.Range("L132").Value = 43521
$ws.Range("M132").Value = -5321
$ws.Range("N132").Value = -48581
$ws.Range("H136").Value = 3679.625
$ws.Range("I136").Value = 2978.739
$ws.Range("J136").Value = 19800
$ws.Range("K136").Value = 8936.217000000001
$ws.Range("L136").Value = 59400
$ws.Range("M136").Value = -6386.217000000001
$ws.Range("N136").Value = -64500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 668.0417
$ws.Range("I7").Value = 906.0769
$ws.Range("K7").Value = 906.0769
$ws.Range("M7").Value = -793.0769
$ws.Range("H31").Value = 26721.348
$ws.Range("J31").Value = 63684.61
$ws.Range("L31").Value = 63684.61
$ws.Range("N31").Value = -64274.61
$ws.Range("H34").Value = 26721.348
$ws.Range("J34").Value = 63684.61
$ws.Range("L34").Value = 63684.61
$ws.Range("N34").Value = -64088.61
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H132").Value = 4314
$ws.Range("I132").Value = 4156.1626
$ws.Range("K132").Value = 12468.4878
$ws.Range("M132").Value = -9938.487799999999
$ws.Range("H135").Value = 69997
$ws.Range("J135").Value = 69997
$ws.Range("L135").Value = 69997
$ws.Range("N135").Value = -80137

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 600
$ws.Range("I59").Value = 600
$ws.Range("K59").Value = 1800
$ws.Range("M59").Value = -1260
$ws.Range("H60").Value = 55555840
$ws.Range("I60").Value = 66666950
$ws.Range("J60").Value = 275
$ws.Range("K60").Value = 200000850
$ws.Range("L60").Value = 825
$ws.Range("M60").Value = -200000599
$ws.Range("N60").Value = -1327
$ws.Range("H61").Value = 978.3570999999999
$ws.Range("I61").Value = 53.8
$ws.Range("J61").Value = 1492
$ws.Range("K61").Value = 161.4
$ws.Range("L61").Value = 4476
$ws.Range("M61").Value = 53.60000000000002
$ws.Range("N61").Value = -4906
$ws.Range("H69").Value = 8213.733
$ws.Range("I69").Value = 5341
$ws.Range("J69").Value = 9650.1
$ws.Range("K69").Value = 16023
$ws.Range("L69").Value = 28950.3
$ws.Range("M69").Value = -15212
$ws.Range("N69").Value = -30572.3
$ws.Range("H72").Value = 8213.733
$ws.Range("I72").Value = 5341
$ws.Range("J72").Value = 9650.1
$ws.Range("K72").Value = 48069
$ws.Range("L72").Value = 86850.90000000001
$ws.Range("M72").Value = -44013
$ws.Range("N72").Value = -94962.90000000001
$ws.Range("H113").Value = 1391.9286
$ws.Range("I113").Value = 856.8
$ws.Range("K113").Value = 2570.4
$ws.Range("M113").Value = -400.3999999999996
$ws.Range("H132").Value = 5663.636
$ws.Range("J132").Value = 6100
$ws.Range("L132").Value = 54900
$ws.Range("N132").Value = -59960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17613.824
$ws.Range("I70").Value = 5403.636
$ws.Range("K70").Value = 5403.636
$ws.Range("M70").Value = -5133.636
$ws.Range("H73").Value = 17613.824
$ws.Range("I73").Value = 5403.636
$ws.Range("K73").Value = 5403.636
$ws.Range("M73").Value = -4467.636

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3469.3845
$ws.Range("J22").Value = 5143.143
$ws.Range("L22").Value = 5143.143
$ws.Range("N22").Value = -5733.143
$ws.Range("H27").Value = 3469.3845
$ws.Range("J27").Value = 5143.143
$ws.Range("L27").Value = 5143.143
$ws.Range("N27").Value = -5357.143
$ws.Range("H46").Value = 4873.875
$ws.Range("I46").Value = 1495
$ws.Range("J46").Value = 6000.1665
$ws.Range("K46").Value = 1495
$ws.Range("L46").Value = 6000.1665
$ws.Range("M46").Value = -1307
$ws.Range("N46").Value = -6376.1665
$ws.Range("H55").Value = 1191954.5
$ws.Range("I55").Value = 1924169.2
$ws.Range("J55").Value = 2105.625
$ws.Range("K55").Value = 1924169.2
$ws.Range("L55").Value = 2105.625
$ws.Range("M55").Value = -1923996.2
$ws.Range("N55").Value = -2451.625
$ws.Range("H68").Value = 6199.625
$ws.Range("I68").Value = 2519
$ws.Range("J68").Value = 12334
$ws.Range("K68").Value = 2519
$ws.Range("L68").Value = 12334
$ws.Range("M68").Value = -1770
$ws.Range("N68").Value = -13832
$ws.Range("H71").Value = 6199.625
$ws.Range("I71").Value = 2519
$ws.Range("J71").Value = 12334
$ws.Range("K71").Value = 12595
$ws.Range("L71").Value = 61670
$ws.Range("M71").Value = -8851
$ws.Range("N71").Value = -69158
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H68").Value = 15000
$ws.Range("J68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16622
$ws.Range("H71").Value = 15000
$ws.Range("J71").Value = 15000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -53112
$ws.Range("H100").Value = 1422.3334
$ws.Range("I100").Value = 607
$ws.Range("J100").Value = 3053
$ws.Range("K100").Value = 1214
$ws.Range("L100").Value = 6106
$ws.Range("M100").Value = -673
$ws.Range("N100").Value = -7188
$ws.Range("H132").Value = 4366.516
$ws.Range("I132").Value = 4155.607
$ws.Range("J132").Value = 6335
$ws.Range("K132").Value = 12466.821
$ws.Range("L132").Value = 19005
$ws.Range("M132").Value = -9936.821
$ws.Range("N132").Value = -24065
